$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 26.83824466666667
$ws.Range("H2").Value = 80.514734
$ws.Range("I2").Value = 0.8882651037973995
$ws.Range("J2").Value = 0.8882651037973996
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.828998666666665
$ws.Range("N2").Value = 26.486996
$ws.Range("O2").Value = 0.1794455804823882
$ws.Range("P2").Value = 0.1794455804823882
$ws.Range("Q2").Value = 236.9548263776738
$ws.Range("R2").Value = 2132.593437399064
$ws.Range("S2").Value = 0.1593952471731731
$ws.Range("T2").Value = 0.1593952471731731

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 26.83824466666667
$ws.Range("H3").Value = 80.514734
$ws.Range("I3").Value = 0.8882651037973995
$ws.Range("J3").Value = 0.8882651037973996
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.33828433333333
$ws.Range("N3").Value = 58.014853
$ws.Range("O3").Value = 0.3930422677296217
$ws.Range("P3").Value = 0.3930422677296217
$ws.Range("Q3").Value = 519.005606371567
$ws.Range("R3").Value = 4671.050457344102
$ws.Range("S3").Value = 0.3491257307416177
$ws.Range("T3").Value = 0.3491257307416177

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 26.83824466666667
$ws.Range("H4").Value = 80.514734
$ws.Range("I4").Value = 0.8882651037973995
$ws.Range("J4").Value = 0.8882651037973996
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.03425566666667
$ws.Range("N4").Value = 63.102767
$ws.Range("O4").Value = 0.4275121517879902
$ws.Range("P4").Value = 0.4275121517879902
$ws.Range("Q4").Value = 564.5224999632198
$ws.Range("R4").Value = 5080.702499668978
$ws.Range("S4").Value = 0.3797441258826087
$ws.Range("T4").Value = 0.3797441258826088

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.239011
$ws.Range("H5").Value = 3.717033
$ws.Range("I5").Value = 0.04100753414354395
$ws.Range("J5").Value = 0.04100753414354396
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.828998666666665
$ws.Range("N5").Value = 26.486996
$ws.Range("O5").Value = 0.1794455804823882
$ws.Range("P5").Value = 0.1794455804823882
$ws.Range("Q5").Value = 10.93922646698533
$ws.Range("R5").Value = 98.45303820286799
$ws.Range("S5").Value = 0.007358620768539596
$ws.Range("T5").Value = 0.007358620768539598

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.239011
$ws.Range("H6").Value = 3.717033
$ws.Range("I6").Value = 0.04100753414354395
$ws.Range("J6").Value = 0.04100753414354396
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.33828433333333
$ws.Range("N6").Value = 58.014853
$ws.Range("O6").Value = 0.3930422677296217
$ws.Range("P6").Value = 0.3930422677296217
$ws.Range("Q6").Value = 23.96034701012767
$ws.Range("R6").Value = 215.643123091149
$ws.Range("S6").Value = 0.01611769421377841
$ws.Range("T6").Value = 0.01611769421377841

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.239011
$ws.Range("H7").Value = 3.717033
$ws.Range("I7").Value = 0.04100753414354395
$ws.Range("J7").Value = 0.04100753414354396
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.03425566666667
$ws.Range("N7").Value = 63.102767
$ws.Range("O7").Value = 0.4275121517879902
$ws.Range("P7").Value = 0.4275121517879902
$ws.Range("Q7").Value = 26.06167414781233
$ws.Range("R7").Value = 234.555067330311
$ws.Range("S7").Value = 0.01753121916122595
$ws.Range("T7").Value = 0.01753121916122596

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.136972666666667
$ws.Range("H8").Value = 6.410918000000001
$ws.Range("I8").Value = 0.07072736205905639
$ws.Range("J8").Value = 0.0707273620590564
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.828998666666665
$ws.Range("N8").Value = 26.486996
$ws.Range("O8").Value = 0.1794455804823882
$ws.Range("P8").Value = 0.1794455804823882
$ws.Range("Q8").Value = 18.86732882470311
$ws.Range("R8").Value = 169.805959422328
$ws.Range("S8").Value = 0.01269171254067541
$ws.Range("T8").Value = 0.01269171254067541

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.136972666666667
$ws.Range("H9").Value = 6.410918000000001
$ws.Range("I9").Value = 0.07072736205905639
$ws.Range("J9").Value = 0.0707273620590564
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.33828433333333
$ws.Range("N9").Value = 58.014853
$ws.Range("O9").Value = 0.3930422677296217
$ws.Range("P9").Value = 0.3930422677296217
$ws.Range("Q9").Value = 41.32538504056156
$ws.Range("R9").Value = 371.928465365054
$ws.Range("S9").Value = 0.02779884277422553
$ws.Range("T9").Value = 0.02779884277422554

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.136972666666667
$ws.Range("H10").Value = 6.410918000000001
$ws.Range("I10").Value = 0.07072736205905639
$ws.Range("J10").Value = 0.0707273620590564
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.03425566666667
$ws.Range("N10").Value = 63.102767
$ws.Range("O10").Value = 0.4275121517879902
$ws.Range("P10").Value = 0.4275121517879902
$ws.Range("Q10").Value = 44.94962942334512
$ws.Range("R10").Value = 404.5466648101061
$ws.Range("S10").Value = 0.03023680674415545
$ws.Range("T10").Value = 0.03023680674415546
